$wb = $excel.ActiveWorkbook

# --- Rename sheets ---------------------------------------------------------
$wsBib = $wb.Worksheets.Item(1)
$wsBib.Name = "Bibliographic Item"

$wsResult = $wb.Worksheets.Item(2)
$wsResult.Name = "Result Set"

$wsProvider = $wb.Worksheets.Item(3)
$wsProvider.Name = "Data Provider"

# --- Populate "Result Set" sheet (sheet2) ---------------------------------
# Column B labels (written top to bottom first so they claim shared-string
# indices 46-50 in order).
$wsResult.Range("B2").Value = "currentPage"
$wsResult.Range("B3").Value = "numberOfResults"
$wsResult.Range("B4").Value = "maxResultsPerPage"
$wsResult.Range("B5").Value = "results"
$wsResult.Range("B6").Value = "warningMessage"

# Column C descriptions (written in this specific order so the resulting
# shared-string indices line up with the target file: 51,52,53,54,55).
$wsResult.Range("C5").Value = "Tableau d'Items"
$wsResult.Range("C4").Value = "Sert à calculer s'il reste des pages de résultats à afficher."
$wsResult.Range("C3").Value = "Sert à mettre à jour les conteneurs ""Stats""."
$wsResult.Range("C2").Value = "Sert à calculer s'il reste des pages de résultats à afficher.`nSert à calculer s'il s'agit du premier ensemble de résultats pour une nouvelle requête."
$wsResult.Range("C6").Value = "Sert à indiquer si la source de données demande une requête moins coûteuse."

$wsResult.Range("C2").WrapText = $true
$wsResult.Rows.Item(2).RowHeight = 43.2

$wsResult.Columns.Item(2).ColumnWidth = 16.88671875
$wsResult.Columns.Item(3).ColumnWidth = 47.33203125

# --- Populate "Data Provider" sheet (sheet3) ------------------------------
$wsProvider.Range("B2").Value = "_BASE_URL"
$wsProvider.Range("B3").Value = "_MAX_RESULTS_PER_PAGE"
$wsProvider.Range("B4").Value = "getSearchResults"
$wsProvider.Range("C4").Value = "searchString, pageNumber"
$wsProvider.Range("B5").Value = "getDetailedItem"
$wsProvider.Range("C5").Value = "url"
$wsProvider.Range("B6").Value = "_buildRequest"
$wsProvider.Range("B7").Value = "_buildResultSet"
$wsProvider.Range("B8").Value = "_buildDataItem"
$wsProvider.Range("B9").Value = "_buildDetailedDataItem"
$wsProvider.Range("C7").Value = "rawXmlData"
$wsProvider.Range("B11").Value = "PROPOSITIONS"
$wsProvider.Range("B12").Value = "getItemById"
$wsProvider.Range("B13").Value = "getNextResults"
$wsProvider.Range("B14").Value = "_currentResultPage"
$wsProvider.Range("B15").Value = "_currentQuery"
$wsProvider.Range("B16").Value = "_currentTotalOfResults"
$wsProvider.Range("C3").Value = "Sert à renseigner le maxResultsPerPage des ResultSets."
$wsProvider.Range("C2").Value = "Servait à calculer la catalogUrl des Items."

# Cells that reuse a string already created above.
$wsProvider.Range("C6").Value = "searchString, pageNumber"
$wsProvider.Range("C8").Value = "rawXmlData"
$wsProvider.Range("C9").Value = "rawXmlData"

$wsProvider.Columns.Item(2).ColumnWidth = 23.5546875
$wsProvider.Columns.Item(3).ColumnWidth = 46.44140625

# --- Selections (must be set in this order: sheet1 keeps its existing
# selection, then sheet2, then sheet3 last so sheet3 ends up active). -----
$null = $wsResult.Range("C7").Select()
$null = $wsProvider.Range("B15").Select()
$wsProvider.Activate()
